# "mobile home page update"
#
# - Mobile sheet: add a new row describing the home page / login-link item,
#   and move the "current cell" down to C11.
# - Website sheet: just the cursor/selection moved on to C41 (no longer the
#   active tab).
# - New "Testing - Revisions" sheet appended at the end of the tab strip,
#   listing the three things that get tested on each revision (Move,
#   Add/delete section, Add/delete clause), and it becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- Website: selection moves to C41 -----------------------------------
$wsWebsite = $wb.Worksheets.Item("Website")
$wsWebsite.Activate()
$wsWebsite.Range("C41").Select()

# --- Mobile: new "home" row + selection moves to C11 --------------------
$wsMobile = $wb.Worksheets.Item("Mobile")
$wsMobile.Activate()
$wsMobile.Range("A7").Value = "A"
$wsMobile.Range("C7").Value = "home"
$wsMobile.Range("D7").Value = "view home page with login link"
$wsMobile.Range("C11").Select()

# --- New sheet: Testing - Revisions, added after the last sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTest = $wb.Worksheets.Add($null, $lastSheet)
$wsTest.Name = "Testing - Revisions"

$wsTest.Columns.Item(1).ColumnWidth = 21.3

$wsTest.Range("A2").Value = "Move"
$wsTest.Range("A4").Value = "Add/delete section"
$wsTest.Range("A3").Value = "Add/delete clause"

$wsTest.Range("B7:C7").Select()
